$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the name for row 7 (ho / ten columns)
$ws.Range("D7").Value = "thanh"
$ws.Range("E7").Value = "phong"

# Update the active selection shown in the sheet view
$ws.Range("G12").Select()
